$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape.
# Price cells are forced to text format first so numeric-looking strings
# (e.g. "2.60", "431.70") keep their exact original textual representation
# instead of being reinterpreted as floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.784.58"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.911.12"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.05"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.58"
$ws.Range("E6").Value = "  -5.86%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.909.55"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.81"
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("E12").Value = "  -4.20%  "
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.47"
$ws.Range("E14").Value = "  -6.10%  "
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.394.05"
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.738.28"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.73"
$ws.Range("E18").Value = "  -5.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.911.25"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.70"
$ws.Range("E20").Value = "  -4.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  -4.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("E23").Value = "  -4.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.25"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.82"
$ws.Range("E25").Value = "  -4.41%  "
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.84"
$ws.Range("E27").Value = "  -4.45%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.48"
$ws.Range("E33").Value = "  -3.99%  "
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0863"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.64"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.00"
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.64"
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.13"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "375.12"
$ws.Range("E45").Value = "  -5.77%  "
$ws.Range("E46").Value = "  -3.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.665.51"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.89"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.26"
$ws.Range("E50").Value = "  -0.84%  "
